$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet ("ODI Bowling")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row (row 1) - copy formatting (bold + border + centered) from an
# existing header cell so we reuse the same style definition.
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

$wb.Worksheets.Item("Player Info").Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Data row (row 2) - MATCH_CODE is stored as text (not a number), matching
# the other sheets in this workbook where "4660" is always text.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4660"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NO"

# Restore the originally active sheet ("Player Info")
$wb.Worksheets.Item("Player Info").Activate()
